# Add a new "% of Q Drop's" column (column I) to the Summer 2015 GE grade
# distribution sheet, mirroring the existing "% of A's".."% of F's" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "% of Q Drop's"

# Per-course "% of Q Drop's" values (stored as text, like the neighboring
# percentage columns, so use a leading apostrophe to stop Excel from
# re-interpreting the text as a numeric percentage).
$qdrop = @{
    3  = "0.00%"
    4  = "1.80%"
    7  = "0.00%"
    10 = "0.00%"
    13 = "0.78%"
    16 = "0.00%"
    19 = "2.70%"
    22 = "0.00%"
    23 = "0.00%"
    26 = "2.08%"
    29 = "0.00%"
    30 = "0.00%"
    31 = "0.00%"
    32 = "0.00%"
    35 = "3.57%"
    38 = "0.00%"
    41 = "3.23%"
    42 = "0.00%"
    45 = "7.69%"
    48 = "0.00%"
}

foreach ($row in $qdrop.Keys) {
    $ws.Cells.Item($row, 9).Value = "'" + $qdrop[$row]
}
